$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R (2021) to the existing year table in row 4 (headers),
# row 5 (percentage series) and row 6 (absolute series), reusing the same
# visual formatting as the neighboring cells so the new column blends in
# with the rest of the table.

# Header cell R4 = 2021, formatted like Q4 (year header style).
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# R5 = 31.8, formatted like D5 (percentage row, "General" style variant).
$ws.Range("D5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 31.8

# R6 = 12957.1, formatted like Q6 (absolute values row).
$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 12957.1

# Update the selection to reflect the newly added column, as in the source
# workbook (active cell R4, selected range R4:R6).
$ws.Range("R4:R6").Select()
